$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relations")
$ws.Range("C1").Value = "src"
$ws.Range("C2").Select() | Out-Null
